# This script applies the day-to-day refresh of the addition/subtraction drill
# sheet: the heading date moves forward one day, and all 100 arithmetic prompts
# in the 20x5 table are swapped for a new set of problems.
#
# Replacements are applied with Find/Replace (Replace = wdReplaceOne, one hit at a
# time) walking the table in its natural reading order (row by row, left to right).
# Doing it one occurrence at a time, in document order, correctly disambiguates the
# sheet's one duplicated prompt, "19+60=" (row 12 col 5 and row 17 col 3), since each
# Find.Execute() call consumes the next remaining occurrence in the document.
#
# One pair is special-cased out of order: row 14 col 4 changes "76-12=" to "78+13=",
# and later, row 20 col 4 changes "8+13=" to "70+29=". Because "8+13=" is a substring
# of "78+13=", replacing row 14 col 4 first would leave a spurious "8+13=" inside it,
# and the row 20 col 4 search would incorrectly match that instead of its own cell. So
# row 20 col 4 is updated first, before it can be shadowed.

$d = $word.ActiveDocument

# Heading date
$d.Content.Find.Execute("2024-10-09 Wednesday", $true, $false, $false, $false, $false, $true, 0, $false, "2024-10-10 Thursday", 1) | Out-Null

# Out-of-order fix-up (see note above) for row 20, col 4, so it is safely updated
# before row 14, col 4 introduces a same-text substring ("8+13=" inside "78+13=").
$d.Content.Find.Execute("8+13=", $true, $false, $false, $false, $false, $true, 0, $false, "70+29=", 1) | Out-Null  # table row 20, col 4

# Remaining 99 table cells, in normal row-by-row, left-to-right order
$d.Content.Find.Execute("7+6=", $true, $false, $false, $false, $false, $true, 0, $false, "91-49=", 1) | Out-Null  # table row 1, col 1
$d.Content.Find.Execute("9+86=", $true, $false, $false, $false, $false, $true, 0, $false, "65-16=", 1) | Out-Null  # table row 1, col 2
$d.Content.Find.Execute("49+42=", $true, $false, $false, $false, $false, $true, 0, $false, "83-77=", 1) | Out-Null  # table row 1, col 3
$d.Content.Find.Execute("40+46=", $true, $false, $false, $false, $false, $true, 0, $false, "74-21=", 1) | Out-Null  # table row 1, col 4
$d.Content.Find.Execute("58-0=", $true, $false, $false, $false, $false, $true, 0, $false, "21+52=", 1) | Out-Null  # table row 1, col 5
$d.Content.Find.Execute("29+41=", $true, $false, $false, $false, $false, $true, 0, $false, "8+81=", 1) | Out-Null  # table row 2, col 1
$d.Content.Find.Execute("54+39=", $true, $false, $false, $false, $false, $true, 0, $false, "20+51=", 1) | Out-Null  # table row 2, col 2
$d.Content.Find.Execute("91-79=", $true, $false, $false, $false, $false, $true, 0, $false, "85-66=", 1) | Out-Null  # table row 2, col 3
$d.Content.Find.Execute("53-48=", $true, $false, $false, $false, $false, $true, 0, $false, "34-14=", 1) | Out-Null  # table row 2, col 4
$d.Content.Find.Execute("49+20=", $true, $false, $false, $false, $false, $true, 0, $false, "53+6=", 1) | Out-Null  # table row 2, col 5
$d.Content.Find.Execute("98-67=", $true, $false, $false, $false, $false, $true, 0, $false, "45-20=", 1) | Out-Null  # table row 3, col 1
$d.Content.Find.Execute("9+15=", $true, $false, $false, $false, $false, $true, 0, $false, "11+32=", 1) | Out-Null  # table row 3, col 2
$d.Content.Find.Execute("64+19=", $true, $false, $false, $false, $false, $true, 0, $false, "3+89=", 1) | Out-Null  # table row 3, col 3
$d.Content.Find.Execute("30+61=", $true, $false, $false, $false, $false, $true, 0, $false, "85-65=", 1) | Out-Null  # table row 3, col 4
$d.Content.Find.Execute("28+44=", $true, $false, $false, $false, $false, $true, 0, $false, "71+17=", 1) | Out-Null  # table row 3, col 5
$d.Content.Find.Execute("78-52=", $true, $false, $false, $false, $false, $true, 0, $false, "84+8=", 1) | Out-Null  # table row 4, col 1
$d.Content.Find.Execute("79-58=", $true, $false, $false, $false, $false, $true, 0, $false, "87-27=", 1) | Out-Null  # table row 4, col 2
$d.Content.Find.Execute("91-57=", $true, $false, $false, $false, $false, $true, 0, $false, "44-22=", 1) | Out-Null  # table row 4, col 3
$d.Content.Find.Execute("57-25=", $true, $false, $false, $false, $false, $true, 0, $false, "90-20=", 1) | Out-Null  # table row 4, col 4
$d.Content.Find.Execute("49-40=", $true, $false, $false, $false, $false, $true, 0, $false, "96-27=", 1) | Out-Null  # table row 4, col 5
$d.Content.Find.Execute("62-27=", $true, $false, $false, $false, $false, $true, 0, $false, "54-48=", 1) | Out-Null  # table row 5, col 1
$d.Content.Find.Execute("39+3=", $true, $false, $false, $false, $false, $true, 0, $false, "2+26=", 1) | Out-Null  # table row 5, col 2
$d.Content.Find.Execute("2+28=", $true, $false, $false, $false, $false, $true, 0, $false, "92-34=", 1) | Out-Null  # table row 5, col 3
$d.Content.Find.Execute("41+44=", $true, $false, $false, $false, $false, $true, 0, $false, "71-47=", 1) | Out-Null  # table row 5, col 4
$d.Content.Find.Execute("84-7=", $true, $false, $false, $false, $false, $true, 0, $false, "19-4=", 1) | Out-Null  # table row 5, col 5
$d.Content.Find.Execute("35+1=", $true, $false, $false, $false, $false, $true, 0, $false, "48+3=", 1) | Out-Null  # table row 6, col 1
$d.Content.Find.Execute("66+1=", $true, $false, $false, $false, $false, $true, 0, $false, "45+18=", 1) | Out-Null  # table row 6, col 2
$d.Content.Find.Execute("47-22=", $true, $false, $false, $false, $false, $true, 0, $false, "10+47=", 1) | Out-Null  # table row 6, col 3
$d.Content.Find.Execute("22+9=", $true, $false, $false, $false, $false, $true, 0, $false, "21+11=", 1) | Out-Null  # table row 6, col 4
$d.Content.Find.Execute("59-23=", $true, $false, $false, $false, $false, $true, 0, $false, "29-1=", 1) | Out-Null  # table row 6, col 5
$d.Content.Find.Execute("71+22=", $true, $false, $false, $false, $false, $true, 0, $false, "14+63=", 1) | Out-Null  # table row 7, col 1
$d.Content.Find.Execute("40+16=", $true, $false, $false, $false, $false, $true, 0, $false, "76-57=", 1) | Out-Null  # table row 7, col 2
$d.Content.Find.Execute("3+43=", $true, $false, $false, $false, $false, $true, 0, $false, "59-51=", 1) | Out-Null  # table row 7, col 3
$d.Content.Find.Execute("26+37=", $true, $false, $false, $false, $false, $true, 0, $false, "27+20=", 1) | Out-Null  # table row 7, col 4
$d.Content.Find.Execute("66+6=", $true, $false, $false, $false, $false, $true, 0, $false, "24+9=", 1) | Out-Null  # table row 7, col 5
$d.Content.Find.Execute("27+37=", $true, $false, $false, $false, $false, $true, 0, $false, "13+36=", 1) | Out-Null  # table row 8, col 1
$d.Content.Find.Execute("51+35=", $true, $false, $false, $false, $false, $true, 0, $false, "44-0=", 1) | Out-Null  # table row 8, col 2
$d.Content.Find.Execute("67-12=", $true, $false, $false, $false, $false, $true, 0, $false, "17+33=", 1) | Out-Null  # table row 8, col 3
$d.Content.Find.Execute("9+54=", $true, $false, $false, $false, $false, $true, 0, $false, "47-40=", 1) | Out-Null  # table row 8, col 4
$d.Content.Find.Execute("33-17=", $true, $false, $false, $false, $false, $true, 0, $false, "5+76=", 1) | Out-Null  # table row 8, col 5
$d.Content.Find.Execute("70-26=", $true, $false, $false, $false, $false, $true, 0, $false, "86-29=", 1) | Out-Null  # table row 9, col 1
$d.Content.Find.Execute("1+49=", $true, $false, $false, $false, $false, $true, 0, $false, "48+1=", 1) | Out-Null  # table row 9, col 2
$d.Content.Find.Execute("22+26=", $true, $false, $false, $false, $false, $true, 0, $false, "2+97=", 1) | Out-Null  # table row 9, col 3
$d.Content.Find.Execute("7+8=", $true, $false, $false, $false, $false, $true, 0, $false, "9+11=", 1) | Out-Null  # table row 9, col 4
$d.Content.Find.Execute("71-11=", $true, $false, $false, $false, $false, $true, 0, $false, "35+51=", 1) | Out-Null  # table row 9, col 5
$d.Content.Find.Execute("39+48=", $true, $false, $false, $false, $false, $true, 0, $false, "75+14=", 1) | Out-Null  # table row 10, col 1
$d.Content.Find.Execute("99-54=", $true, $false, $false, $false, $false, $true, 0, $false, "13+66=", 1) | Out-Null  # table row 10, col 2
$d.Content.Find.Execute("15+63=", $true, $false, $false, $false, $false, $true, 0, $false, "57+24=", 1) | Out-Null  # table row 10, col 3
$d.Content.Find.Execute("82+17=", $true, $false, $false, $false, $false, $true, 0, $false, "60+6=", 1) | Out-Null  # table row 10, col 4
$d.Content.Find.Execute("3+55=", $true, $false, $false, $false, $false, $true, 0, $false, "67-11=", 1) | Out-Null  # table row 10, col 5
$d.Content.Find.Execute("55-6=", $true, $false, $false, $false, $false, $true, 0, $false, "78-20=", 1) | Out-Null  # table row 11, col 1
$d.Content.Find.Execute("20+33=", $true, $false, $false, $false, $false, $true, 0, $false, "45+22=", 1) | Out-Null  # table row 11, col 2
$d.Content.Find.Execute("38+20=", $true, $false, $false, $false, $false, $true, 0, $false, "69+6=", 1) | Out-Null  # table row 11, col 3
$d.Content.Find.Execute("89-39=", $true, $false, $false, $false, $false, $true, 0, $false, "34-5=", 1) | Out-Null  # table row 11, col 4
$d.Content.Find.Execute("38-13=", $true, $false, $false, $false, $false, $true, 0, $false, "76-13=", 1) | Out-Null  # table row 11, col 5
$d.Content.Find.Execute("59+39=", $true, $false, $false, $false, $false, $true, 0, $false, "78-32=", 1) | Out-Null  # table row 12, col 1
$d.Content.Find.Execute("55+29=", $true, $false, $false, $false, $false, $true, 0, $false, "47-46=", 1) | Out-Null  # table row 12, col 2
$d.Content.Find.Execute("55+44=", $true, $false, $false, $false, $false, $true, 0, $false, "78-39=", 1) | Out-Null  # table row 12, col 3
$d.Content.Find.Execute("9+2=", $true, $false, $false, $false, $false, $true, 0, $false, "23+8=", 1) | Out-Null  # table row 12, col 4
$d.Content.Find.Execute("19+60=", $true, $false, $false, $false, $false, $true, 0, $false, "89-69=", 1) | Out-Null  # table row 12, col 5
$d.Content.Find.Execute("3+16=", $true, $false, $false, $false, $false, $true, 0, $false, "86-5=", 1) | Out-Null  # table row 13, col 1
$d.Content.Find.Execute("10+1=", $true, $false, $false, $false, $false, $true, 0, $false, "21-3=", 1) | Out-Null  # table row 13, col 2
$d.Content.Find.Execute("74-0=", $true, $false, $false, $false, $false, $true, 0, $false, "56-36=", 1) | Out-Null  # table row 13, col 3
$d.Content.Find.Execute("22+55=", $true, $false, $false, $false, $false, $true, 0, $false, "35-31=", 1) | Out-Null  # table row 13, col 4
$d.Content.Find.Execute("9+27=", $true, $false, $false, $false, $false, $true, 0, $false, "54-17=", 1) | Out-Null  # table row 13, col 5
$d.Content.Find.Execute("87-22=", $true, $false, $false, $false, $false, $true, 0, $false, "3+33=", 1) | Out-Null  # table row 14, col 1
$d.Content.Find.Execute("86-11=", $true, $false, $false, $false, $false, $true, 0, $false, "62-10=", 1) | Out-Null  # table row 14, col 2
$d.Content.Find.Execute("70-60=", $true, $false, $false, $false, $false, $true, 0, $false, "84-50=", 1) | Out-Null  # table row 14, col 3
$d.Content.Find.Execute("76-12=", $true, $false, $false, $false, $false, $true, 0, $false, "78+13=", 1) | Out-Null  # table row 14, col 4
$d.Content.Find.Execute("0+79=", $true, $false, $false, $false, $false, $true, 0, $false, "2+61=", 1) | Out-Null  # table row 14, col 5
$d.Content.Find.Execute("60-13=", $true, $false, $false, $false, $false, $true, 0, $false, "32+8=", 1) | Out-Null  # table row 15, col 1
$d.Content.Find.Execute("13+67=", $true, $false, $false, $false, $false, $true, 0, $false, "26+23=", 1) | Out-Null  # table row 15, col 2
$d.Content.Find.Execute("91-21=", $true, $false, $false, $false, $false, $true, 0, $false, "94-1=", 1) | Out-Null  # table row 15, col 3
$d.Content.Find.Execute("64-6=", $true, $false, $false, $false, $false, $true, 0, $false, "78-34=", 1) | Out-Null  # table row 15, col 4
$d.Content.Find.Execute("99-1=", $true, $false, $false, $false, $false, $true, 0, $false, "92-66=", 1) | Out-Null  # table row 15, col 5
$d.Content.Find.Execute("72-34=", $true, $false, $false, $false, $false, $true, 0, $false, "35+6=", 1) | Out-Null  # table row 16, col 1
$d.Content.Find.Execute("40-33=", $true, $false, $false, $false, $false, $true, 0, $false, "2+22=", 1) | Out-Null  # table row 16, col 2
$d.Content.Find.Execute("65-28=", $true, $false, $false, $false, $false, $true, 0, $false, "36-18=", 1) | Out-Null  # table row 16, col 3
$d.Content.Find.Execute("37+44=", $true, $false, $false, $false, $false, $true, 0, $false, "47+18=", 1) | Out-Null  # table row 16, col 4
$d.Content.Find.Execute("40-1=", $true, $false, $false, $false, $false, $true, 0, $false, "60-45=", 1) | Out-Null  # table row 16, col 5
$d.Content.Find.Execute("15+58=", $true, $false, $false, $false, $false, $true, 0, $false, "90-11=", 1) | Out-Null  # table row 17, col 1
$d.Content.Find.Execute("8+19=", $true, $false, $false, $false, $false, $true, 0, $false, "36+24=", 1) | Out-Null  # table row 17, col 2
$d.Content.Find.Execute("19+60=", $true, $false, $false, $false, $false, $true, 0, $false, "0+43=", 1) | Out-Null  # table row 17, col 3
$d.Content.Find.Execute("51+26=", $true, $false, $false, $false, $false, $true, 0, $false, "71-8=", 1) | Out-Null  # table row 17, col 4
$d.Content.Find.Execute("24+35=", $true, $false, $false, $false, $false, $true, 0, $false, "25+56=", 1) | Out-Null  # table row 17, col 5
$d.Content.Find.Execute("75-49=", $true, $false, $false, $false, $false, $true, 0, $false, "95-87=", 1) | Out-Null  # table row 18, col 1
$d.Content.Find.Execute("87-74=", $true, $false, $false, $false, $false, $true, 0, $false, "76+21=", 1) | Out-Null  # table row 18, col 2
$d.Content.Find.Execute("16+1=", $true, $false, $false, $false, $false, $true, 0, $false, "17+11=", 1) | Out-Null  # table row 18, col 3
$d.Content.Find.Execute("29+38=", $true, $false, $false, $false, $false, $true, 0, $false, "1+83=", 1) | Out-Null  # table row 18, col 4
$d.Content.Find.Execute("32+10=", $true, $false, $false, $false, $false, $true, 0, $false, "84-35=", 1) | Out-Null  # table row 18, col 5
$d.Content.Find.Execute("82-44=", $true, $false, $false, $false, $false, $true, 0, $false, "58+5=", 1) | Out-Null  # table row 19, col 1
$d.Content.Find.Execute("14+67=", $true, $false, $false, $false, $false, $true, 0, $false, "75-27=", 1) | Out-Null  # table row 19, col 2
$d.Content.Find.Execute("22+64=", $true, $false, $false, $false, $false, $true, 0, $false, "24+56=", 1) | Out-Null  # table row 19, col 3
$d.Content.Find.Execute("57-54=", $true, $false, $false, $false, $false, $true, 0, $false, "34-20=", 1) | Out-Null  # table row 19, col 4
$d.Content.Find.Execute("59-2=", $true, $false, $false, $false, $false, $true, 0, $false, "61-18=", 1) | Out-Null  # table row 19, col 5
$d.Content.Find.Execute("99-50=", $true, $false, $false, $false, $false, $true, 0, $false, "26+58=", 1) | Out-Null  # table row 20, col 1
$d.Content.Find.Execute("10+6=", $true, $false, $false, $false, $false, $true, 0, $false, "61-48=", 1) | Out-Null  # table row 20, col 2
$d.Content.Find.Execute("69-47=", $true, $false, $false, $false, $false, $true, 0, $false, "37+20=", 1) | Out-Null  # table row 20, col 3
$d.Content.Find.Execute("58+14=", $true, $false, $false, $false, $false, $true, 0, $false, "72-40=", 1) | Out-Null  # table row 20, col 5
